$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.718.94"
$ws.Range("E2").Value = "  +6.11%  "

$ws.Range("D3").Value = "3.006.58"
$ws.Range("E3").Value = "  +3.45%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.18"
$ws.Range("E5").Value = "  +2.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.56"
$ws.Range("E6").Value = "  +12.14%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.518"
$ws.Range("E8").Value = "  +3.40%  "

$ws.Range("D9").Value = "3.002.17"
$ws.Range("E9").Value = "  +3.39%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.71"
$ws.Range("E10").Value = "  -4.14%  "

$ws.Range("E11").Value = "  +6.67%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.459"
$ws.Range("E12").Value = "  +7.26%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000250"
$ws.Range("E13").Value = "  +8.36%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.66"
$ws.Range("E14").Value = "  +7.54%  "

$ws.Range("E15").Value = "  -0.68%  "

$ws.Range("D16").Value = "65.770.11"
$ws.Range("E16").Value = "  +6.27%  "

$ws.Range("D17").Value = "3.506.77"
$ws.Range("E17").Value = "  +3.51%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.96"
$ws.Range("E18").Value = "  +7.11%  "

$ws.Range("D19").Value = "3.008.57"
$ws.Range("E19").Value = "  +3.68%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "456.96"
$ws.Range("E20").Value = "  +6.46%  "

$ws.Range("E21").Value = "  +7.89%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.688"
$ws.Range("E22").Value = "  +5.62%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.36"
$ws.Range("E23").Value = "  +7.24%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.36"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.27"
$ws.Range("E25").Value = "  +12.33%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.34"
$ws.Range("E26").Value = "  +2.50%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.67"
$ws.Range("E27").Value = "  +5.09%  "

$ws.Range("E28").Value = "  -0.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.06"
$ws.Range("E29").Value = "  +15.36%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.35"
$ws.Range("E30").Value = "  +16.68%  "

$ws.Range("E31").Value = "  -5.59%  "

$ws.Range("E32").Value = "  +3.46%  "

$ws.Range("E33").Value = "  +5.50%  "

$ws.Range("E34").Value = "  +3.81%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.989"
$ws.Range("E36").Value = "  +3.50%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.81"
$ws.Range("E37").Value = "  +8.09%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.15"
$ws.Range("E38").Value = "  +12.96%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.70"
$ws.Range("E39").Value = "  +1.83%  "

$ws.Range("E40").Value = "  +2.68%  "

$ws.Range("E41").Value = "  +16.34%  "

$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.122"
$ws.Range("E42").Value = "  +6.43%  "

$ws.Range("B43").Value = "Arweave"
$ws.Range("C43").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "43.33"
$ws.Range("E43").Value = "  +5.33%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.45"
$ws.Range("E44").Value = "  +3.60%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "392.11"
$ws.Range("E45").Value = "  +13.10%  "

$ws.Range("D46").Value = "2.789.59"
$ws.Range("E46").Value = "  +3.24%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0353"
$ws.Range("E47").Value = "  +5.32%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.90"
$ws.Range("E48").Value = "  +2.31%  "

$ws.Range("E49").Value = "  -0.08%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.71"
$ws.Range("E50").Value = "  +9.92%  "

$ws.Range("E51").Value = "  +4.15%  "
